$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 17, column D: status changes from "实审" to "受理"
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = "受理"

# ---------------------------------------------------------------------------
# 2) New row 19 & 20 content first borrows formatting from the (still blank)
#    placeholder row 18 (style "5" cells) and from the header row (style "1"
#    cells) and from row 2 (style "6" cell) -- do this BEFORE row 18 itself
#    is overwritten with its own new content/style below.
# ---------------------------------------------------------------------------

# Row 19 A/B pick up the same format the blank placeholder row (old row 18)
# had (font2, vertical-center + wrap, no explicit horizontal alignment).
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial(-4122)

# Row 19 C/D reuse the "center, wrap" style already used throughout column C/D.
$ws.Range("C17").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("D19").PasteSpecial(-4122)

# Row 19 E: start from the blank placeholder style, then force a text number
# format on it -- this mirrors a new, unique cell style (wrap + vertical
# center, no horizontal override, numeric format "@").
$ws.Range("E18").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").NumberFormat = "@"

# Row 20 A/B reuse the header-row style (font0, vertical-center + wrap).
$ws.Range("A1").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("B20").PasteSpecial(-4122)

# Row 20 C/D reuse the "center, wrap" style.
$ws.Range("C17").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("D20").PasteSpecial(-4122)

# Row 20 E reuses the style already used for column E elsewhere (row 2).
$ws.Range("E2").Copy()
$ws.Range("E20").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Row 18 becomes a real data row; give it the same look as row 17.
# ---------------------------------------------------------------------------
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Fill in the cell values for the three new patent/publication rows.
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "刘江,李衡,俞向阳,区明阳,李昊锦,胡衍"
$ws.Range("B18").Value = "医学图像生成方法、装置、电子设备及存储介质"
$ws.Range("C18").Value = "发明"
$ws.Range("D18").Value = "受理"
$ws.Range("E18").Value = "2023108192704"

$ws.Range("A19").Value = "刘江,胡衍,沈俊勇"
$ws.Range("B19").Value = "病变检测模型的训练方法和装置、电子设备及存储介质"
$ws.Range("C19").Value = "发明"
$ws.Range("D19").Value = "受理"
$ws.Range("E19").Value = "2023105094659"

$ws.Range("A20").Value = "刘江,胡衍,叶海礼,陈晓慧"
$ws.Range("B20").Value = "眼底图像预测方法、眼底图像预测系统、设备及存储介质"
$ws.Range("C20").Value = "发明"
$ws.Range("D20").Value = "实审"
$ws.Range("E20").Value = "2023102311682"

# ---------------------------------------------------------------------------
# 5) Row heights: rows with wrapped two-line text render at 42pt.
# ---------------------------------------------------------------------------
$ws.Rows("18:18").RowHeight = 42
$ws.Rows("19:19").RowHeight = 42
$ws.Rows("20:20").RowHeight = 42

# ---------------------------------------------------------------------------
# 6) Update the visible selection to the newly-added last row.
# ---------------------------------------------------------------------------
$ws.Range("A18:E18").Select()
